# ============================================================
# Add "2022-Q3" sheet as the new second sheet (after "总计"),
# shifting the existing quarterly sheets down by one position.
# Populate it with the fund-holding detail table, and update
# the "总计" (total) summary sheet with a new row for 2022-Q3.
# ============================================================

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2    = $wb.Worksheets.Item(2)   # currently "2022-Q2", will stay in place; new sheet inserted before it

# --- 1. Insert the new worksheet right before the current "2022-Q2" sheet ---
$ws2 = $wb.Worksheets.Add($wsQ2)
$ws2.Name = "2022-Q3"

# --- 2. Populate '2022-Q3' worksheet header + data ---
$ws2.Cells.Item(1,2).Value = "基金代码"
$ws2.Cells.Item(1,3).Value = "基金名称"
$ws2.Cells.Item(1,4).Value = "基金规模"
$ws2.Cells.Item(1,5).Value = "股票总仓位"
$ws2.Cells.Item(1,6).Value = "仓位占比"
$ws2.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws2.Cells.Item(1,8).Value = "仓位排名"
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,2).Value = "'001576"
$ws2.Cells.Item(2,3).Value = "国泰智能装备股票A"
$ws2.Cells.Item(2,4).Value = "'22.66"
$ws2.Cells.Item(2,5).Value = "'93.43"
$ws2.Cells.Item(2,6).Value = "'6.35"
$ws2.Cells.Item(2,7).Value = "'1.4389"
$ws2.Cells.Item(2,8).Value = 6
$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,2).Value = "'340008"
$ws2.Cells.Item(3,3).Value = "兴全有机增长混合"
$ws2.Cells.Item(3,4).Value = "'23.13"
$ws2.Cells.Item(3,5).Value = "'76.53"
$ws2.Cells.Item(3,6).Value = "'3.27"
$ws2.Cells.Item(3,7).Value = "'0.7564"
$ws2.Cells.Item(3,8).Value = 5
$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,2).Value = "'010330"
$ws2.Cells.Item(4,3).Value = "东吴兴享成长混合A"
$ws2.Cells.Item(4,4).Value = "'10.46"
$ws2.Cells.Item(4,5).Value = "'73.70"
$ws2.Cells.Item(4,6).Value = "'6.42"
$ws2.Cells.Item(4,7).Value = "'0.6715"
$ws2.Cells.Item(4,8).Value = 2
$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,2).Value = "'011322"
$ws2.Cells.Item(5,3).Value = "国泰智能装备股票C"
$ws2.Cells.Item(5,4).Value = "'3.96"
$ws2.Cells.Item(5,5).Value = "'93.43"
$ws2.Cells.Item(5,6).Value = "'6.35"
$ws2.Cells.Item(5,7).Value = "'0.2515"
$ws2.Cells.Item(5,8).Value = 6
$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,2).Value = "'011462"
$ws2.Cells.Item(6,3).Value = "东吴兴享成长混合C"
$ws2.Cells.Item(6,4).Value = "'0.81"
$ws2.Cells.Item(6,5).Value = "'73.70"
$ws2.Cells.Item(6,6).Value = "'6.42"
$ws2.Cells.Item(6,7).Value = "'0.0520"
$ws2.Cells.Item(6,8).Value = 2
$ws2.Cells.Item(7,1).Value = 5
$ws2.Cells.Item(7,2).Value = "'005413"
$ws2.Cells.Item(7,3).Value = "金信民长灵活配置混合C"
$ws2.Cells.Item(7,4).Value = "'0.91"
$ws2.Cells.Item(7,5).Value = "'89.47"
$ws2.Cells.Item(7,6).Value = "'4.94"
$ws2.Cells.Item(7,7).Value = "'0.0450"
$ws2.Cells.Item(7,8).Value = 6
$ws2.Cells.Item(8,1).Value = 6
$ws2.Cells.Item(8,2).Value = "'005412"
$ws2.Cells.Item(8,3).Value = "金信民长灵活配置混合A"
$ws2.Cells.Item(8,4).Value = "'0.86"
$ws2.Cells.Item(8,5).Value = "'89.47"
$ws2.Cells.Item(8,6).Value = "'4.94"
$ws2.Cells.Item(8,7).Value = "'0.0425"
$ws2.Cells.Item(8,8).Value = 6
$ws2.Cells.Item(9,1).Value = 7
$ws2.Cells.Item(9,2).Value = "'015694"
$ws2.Cells.Item(9,3).Value = "瑞达策略优选混合A"
$ws2.Cells.Item(9,4).Value = "'0.09"
$ws2.Cells.Item(9,5).Value = "'67.87"
$ws2.Cells.Item(9,6).Value = "'2.66"
$ws2.Cells.Item(9,7).Value = "'0.0024"
$ws2.Cells.Item(9,8).Value = 9
$ws2.Cells.Item(10,1).Value = 8
$ws2.Cells.Item(10,2).Value = "'015695"
$ws2.Cells.Item(10,3).Value = "瑞达策略优选混合C"
$ws2.Cells.Item(10,4).Value = "'0.00"
$ws2.Cells.Item(10,5).Value = "'67.87"
$ws2.Cells.Item(10,6).Value = "'2.66"
$ws2.Cells.Item(10,7).Value = 0
$ws2.Cells.Item(10,8).Value = 9

# --- 3. Apply header/index styling (bold, centered, bordered) to row 1 (B:H) and column A (2:10) ---
$wsTotal.Cells.Item(1,2).Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Cells.Item(1,2).Copy()
$ws2.Range("A2:A10").PasteSpecial(-4122)

# --- 4. Update "总计" sheet: insert new 2022-Q3 row, shift the rest down ---
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $wsTotal.Cells.Item($dest, 2).Value = $wsTotal.Cells.Item($r, 2).Value2
    $wsTotal.Cells.Item($dest, 3).Value = $wsTotal.Cells.Item($r, 3).Value2
    $wsTotal.Cells.Item($dest, 4).Value = $wsTotal.Cells.Item($r, 4).Value2
}

# Preserve the index-column style for the newly exposed last row (row 9)
$wsTotal.Cells.Item(8,1).Copy()
$wsTotal.Cells.Item(9,1).PasteSpecial(-4122)

# New 2022-Q3 figures
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 9
$wsTotal.Cells.Item(2,4).Value = 3.26

# Refresh the running index numbers in column A (0..7)
for ($r = 2; $r -le 9; $r++) {
    $wsTotal.Cells.Item($r,1).Value = $r - 2
}
